# Auto-generated Excel COM-interop script to apply the diff changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 2383
$ws.Cells.Item(3, 6).Value = 589
$ws.Cells.Item(4, 6).Value = 213
$ws.Cells.Item(5, 6).Value = 370
$ws.Cells.Item(6, 6).Value = 370
$ws.Cells.Item(7, 6).Value = 623
$ws.Cells.Item(9, 6).Value = 830
$ws.Cells.Item(10, 6).Value = 547
$ws.Cells.Item(11, 6).Value = 859
$ws.Cells.Item(12, 6).Value = 394
$ws.Cells.Item(13, 6).Value = 104
$ws.Cells.Item(14, 6).Value = 410
$ws.Cells.Item(15, 6).Value = 26
$ws.Cells.Item(16, 6).Value = 1044
$ws.Cells.Item(17, 6).Value = 22176
$ws.Cells.Item(18, 2).Value = '2024-07-19'
$ws.Cells.Item(18, 3).Value = '广州·萤火虫动漫游戏嘉年华 × KKWORLD2024 快看漫画乐园'
$ws.Cells.Item(18, 4).Value = '新港东路1000号 保利世贸博览馆'
$ws.Cells.Item(18, 5).Value = '2024.07.19 09:00-07.22 17:00'
$ws.Cells.Item(18, 6).Value = 22177
$ws.Cells.Item(18, 7).Value = '已售罄'
$ws.Cells.Item(18, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87210'
$ws.Cells.Item(18, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/DTCdOTPs1718177177472.jpeg'
$ws.Cells.Item(19, 3).Value = '广州·AP动漫游戏嘉年华'
$ws.Cells.Item(19, 4).Value = '新港东路630-638号 南丰国际会展中心'
$ws.Cells.Item(19, 5).Value = '2024.07.27 09:00-07.28 17:00'
$ws.Cells.Item(19, 6).Value = 1040
$ws.Cells.Item(19, 7).Value = 80
$ws.Cells.Item(19, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87213'
$ws.Cells.Item(19, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/3Z8rGZPP1718164976101.jpeg'
$ws.Cells.Item(20, 3).Value = '广州·LookLook动漫嘉年华2th'
$ws.Cells.Item(20, 4).Value = '东沙大道16号 健康方舟6层博览馆'
$ws.Cells.Item(20, 5).Value = '2024.07.27 10:00-07.28 17:30'
$ws.Cells.Item(20, 6).Value = 102
$ws.Cells.Item(20, 7).Value = 29.9
$ws.Cells.Item(20, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87217'
$ws.Cells.Item(20, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/p4Bs2soo1718113055492.png'
$ws.Cells.Item(21, 2).Value = '2024-07-27'
$ws.Cells.Item(21, 3).Value = '广州·原神x星穹x崩only'
$ws.Cells.Item(21, 4).Value = '鸿盛二路巨大创意产业园 巨大产业园·智汇港'
$ws.Cells.Item(21, 5).Value = '2024.07.27 10:00-07.27 17:00'
$ws.Cells.Item(21, 6).Value = 295
$ws.Cells.Item(21, 7).Value = 55
$ws.Cells.Item(21, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87184'
$ws.Cells.Item(21, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/u67hjpFi1718160712051.jpeg'
$ws.Cells.Item(22, 6).Value = 184
$ws.Cells.Item(23, 6).Value = 187
$ws.Cells.Item(24, 6).Value = 15
$ws.Cells.Item(25, 6).Value = 24
$ws.Cells.Item(26, 6).Value = 272
$ws.Cells.Item(27, 6).Value = 20
$ws.Cells.Item(28, 6).Value = 377
$ws.Cells.Item(29, 6).Value = 167
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(6, 6).Value = 213
$ws.Cells.Item(7, 6).Value = 235
$ws.Cells.Item(8, 6).Value = 3481
$ws.Cells.Item(10, 6).Value = 124
$ws.Cells.Item(16, 6).Value = 4028
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 274
$ws.Cells.Item(3, 6).Value = 128
$ws.Cells.Item(4, 6).Value = 661
$ws.Cells.Item(5, 6).Value = 219
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 274
$ws.Cells.Item(3, 6).Value = 128
$ws.Cells.Item(5, 6).Value = 2383
$ws.Cells.Item(6, 6).Value = 661
$ws.Cells.Item(7, 6).Value = 589
$ws.Cells.Item(8, 6).Value = 213
$ws.Cells.Item(9, 6).Value = 370
$ws.Cells.Item(10, 6).Value = 370
$ws.Cells.Item(11, 6).Value = 623
$ws.Cells.Item(16, 6).Value = 213
$ws.Cells.Item(17, 6).Value = 219
$ws.Cells.Item(18, 6).Value = 830
$ws.Cells.Item(19, 6).Value = 547
$ws.Cells.Item(20, 6).Value = 859
$ws.Cells.Item(21, 6).Value = 394
$ws.Cells.Item(22, 6).Value = 104
$ws.Cells.Item(23, 6).Value = 410
$ws.Cells.Item(24, 3).Value = '广州·火影only'
$ws.Cells.Item(24, 4).Value = '人和镇蚌湖清河大街168号 人和园'
$ws.Cells.Item(24, 5).Value = '2024.07.14 09:30-07.14 17:30'
$ws.Cells.Item(24, 6).Value = 1044
$ws.Cells.Item(24, 7).Value = 78
$ws.Cells.Item(24, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84815'
$ws.Cells.Item(24, 9).Value = '//i2.hdslb.com/bfs/openplatform/202404/QLOhW4Nr1714384036670.png'
$ws.Cells.Item(25, 2).Value = '2024-07-19'
$ws.Cells.Item(25, 3).Value = '广州·萤火虫动漫游戏嘉年华 × KKWORLD2024 快看漫画乐园'
$ws.Cells.Item(25, 4).Value = '新港东路1000号 保利世贸博览馆'
$ws.Cells.Item(25, 5).Value = '2024.07.19 09:00-07.22 17:00'
$ws.Cells.Item(25, 6).Value = 22177
$ws.Cells.Item(25, 7).Value = '已售罄'
$ws.Cells.Item(25, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87210'
$ws.Cells.Item(25, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/DTCdOTPs1718177177472.jpeg'
$ws.Cells.Item(26, 2).Value = '2024-07-20'
$ws.Cells.Item(26, 3).Value = '广州·跨越二次元ACG神级动漫世界巡回演唱会'
$ws.Cells.Item(26, 4).Value = '广州市荔湾区十甫路125号(上下九步行街内)2层 广州平安大戏院'
$ws.Cells.Item(26, 5).Value = '2024.07.20 19:30-07.20 21:10'
$ws.Cells.Item(26, 6).Value = 235
$ws.Cells.Item(26, 7).Value = 280
$ws.Cells.Item(26, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=85353'
$ws.Cells.Item(26, 9).Value = '//i1.hdslb.com/bfs/openplatform/202405/4gACWbPh1715223804704.jpeg'
$ws.Cells.Item(27, 6).Value = 235
$ws.Cells.Item(28, 6).Value = 3481
$ws.Cells.Item(30, 6).Value = 124
$ws.Cells.Item(32, 6).Value = 1040
$ws.Cells.Item(33, 6).Value = 102
$ws.Cells.Item(34, 6).Value = 295
$ws.Cells.Item(37, 6).Value = 322
$ws.Cells.Item(38, 6).Value = 184
$ws.Cells.Item(39, 6).Value = 187
$ws.Cells.Item(40, 6).Value = 15
$ws.Cells.Item(41, 6).Value = 24
$ws.Cells.Item(44, 6).Value = 272
$ws.Cells.Item(45, 6).Value = 20
$ws.Cells.Item(46, 6).Value = 377
$ws.Cells.Item(47, 6).Value = 167
$ws.Cells.Item(48, 6).Value = 4028
